$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "281474991109254-1748625109376"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "Mobile Usage"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2025-05-30T11:11:49.376"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "281474991109254"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "117"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "51834027"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "JORGE HUMBERTO TORNERO"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 20.72264935
$ws.Range("I2").Value = -103.48967802
$ws.Range("J2").Value = 0
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109254/1748625106876/t8BT6yo6Du-camera-video-segment-driver-1748625109376.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAJAEKMH6%2F20250531%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250531T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCuUK7w6mmYqW2%2FcaY5Zsk%2FUjw3OSnuIDDgNxk%2FJCULYAIgfudyuGCLrAYHoTeAmUhjEBGJ6BPZjUG13vKGt9tCUXQq5gMItv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDF02x%2Bud36UupgbyJiq6Aw9aUI9rm6kXCltcBRlv%2FdDMFxt99agNdK3vi4%2B%2BQZuAKhCALjB2TN2EZA6efhNOAW4y0tORrNVnoWdOmNHXFEh6JlKS6dv1DKPdRav05xoccU4gpLASbjizN2RzhSurAgyato%2FfK2fpMWm1sj1O%2B1rZs8K11fG0Qkr5f7N9emLDENT4xJv%2Fi7m9IKd%2BLf7xHZTfWk8xLaHes4ncg2hO%2FZeMXjQcjum7UQHyWcS2Rfhb7kzA74fxHtjfSW3MtjSU31Hi2NdJVJCNHs8ZX%2FIFZmPJH9gq5hWj3NXAFjpevz5zz0ZWovkGKq7D00RRfp9GqIK%2BYSzYNJ8XqJACflqgQtDzS9ymXvbMK2FPCzNG1667r%2FR6%2BJPSUupGrPwkRw%2B4QmxApO6HX1NzegKE8%2BlUXvyAMGEuyVoVj4uokf5IP1P9%2FgWAlxSAbfRwGhB5VfG9a8rYERYnO8IDU%2B9Tta2MlusSEp5siXPq2PR8QBMqoC%2F8zM7gM7BwND6ZILuXNiqaQbGD3XJ96IYGzAutb1mVofPzyjcBzolV2T88xKwv0FGrE%2FpJPzJEPMl5ouZ6oBxxTDiQO10x53lFQJ0w3pvqwQY6pQGoJi3k6bpfI1as9VEqVsSHHieB00CrBHFBlN%2BX7K7zlpMS12ZGN%2BbSoHkhsk2zLCzRsmc%2BShOln5hQBKyLJfDHn01tz0wDxXnXFyrkzVM4kWjVvA4yMXdN3D%2B%2BuSX%2FcY5bQPijbW78%2FRKP82bnDhqn4yygMiBgsM9j%2BIKWOWaw6mkZpNdTRDzFmRJDEzAj94RO4QEUTNpEd%2B8%2FvOHlpFzkvooAGrg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2031%20May%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=3fc508100835727ac96a9dbf0a85588044831f9e3d32a867732d34b7a8f4280d"
$ws.Range("K2").Style = "Normal"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991109254/1748625106876/XEZPQddM7k-camera-video-segment-1748625109376.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAJAEKMH6%2F20250531%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250531T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCuUK7w6mmYqW2%2FcaY5Zsk%2FUjw3OSnuIDDgNxk%2FJCULYAIgfudyuGCLrAYHoTeAmUhjEBGJ6BPZjUG13vKGt9tCUXQq5gMItv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDF02x%2Bud36UupgbyJiq6Aw9aUI9rm6kXCltcBRlv%2FdDMFxt99agNdK3vi4%2B%2BQZuAKhCALjB2TN2EZA6efhNOAW4y0tORrNVnoWdOmNHXFEh6JlKS6dv1DKPdRav05xoccU4gpLASbjizN2RzhSurAgyato%2FfK2fpMWm1sj1O%2B1rZs8K11fG0Qkr5f7N9emLDENT4xJv%2Fi7m9IKd%2BLf7xHZTfWk8xLaHes4ncg2hO%2FZeMXjQcjum7UQHyWcS2Rfhb7kzA74fxHtjfSW3MtjSU31Hi2NdJVJCNHs8ZX%2FIFZmPJH9gq5hWj3NXAFjpevz5zz0ZWovkGKq7D00RRfp9GqIK%2BYSzYNJ8XqJACflqgQtDzS9ymXvbMK2FPCzNG1667r%2FR6%2BJPSUupGrPwkRw%2B4QmxApO6HX1NzegKE8%2BlUXvyAMGEuyVoVj4uokf5IP1P9%2FgWAlxSAbfRwGhB5VfG9a8rYERYnO8IDU%2B9Tta2MlusSEp5siXPq2PR8QBMqoC%2F8zM7gM7BwND6ZILuXNiqaQbGD3XJ96IYGzAutb1mVofPzyjcBzolV2T88xKwv0FGrE%2FpJPzJEPMl5ouZ6oBxxTDiQO10x53lFQJ0w3pvqwQY6pQGoJi3k6bpfI1as9VEqVsSHHieB00CrBHFBlN%2BX7K7zlpMS12ZGN%2BbSoHkhsk2zLCzRsmc%2BShOln5hQBKyLJfDHn01tz0wDxXnXFyrkzVM4kWjVvA4yMXdN3D%2B%2BuSX%2FcY5bQPijbW78%2FRKP82bnDhqn4yygMiBgsM9j%2BIKWOWaw6mkZpNdTRDzFmRJDEzAj94RO4QEUTNpEd%2B8%2FvOHlpFzkvooAGrg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2031%20May%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=c7c4beffef35b86a69c4f69a937aa95d28e6ba83b4dc7ae64651ef900e0b7e9e"
$ws.Range("L2").Style = "Normal"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "281474991109254-1748621868279"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "Mobile Usage"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2025-05-30T10:17:48.279"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "281474991109254"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "117"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "51834027"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "JORGE HUMBERTO TORNERO"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 20.68034898
$ws.Range("I3").Value = -103.46086208
$ws.Range("J3").Value = 0
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109254/1748621865779/zpEYZ7igT4-camera-video-segment-driver-1748621868279.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAJAEKMH6%2F20250531%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250531T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCuUK7w6mmYqW2%2FcaY5Zsk%2FUjw3OSnuIDDgNxk%2FJCULYAIgfudyuGCLrAYHoTeAmUhjEBGJ6BPZjUG13vKGt9tCUXQq5gMItv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDF02x%2Bud36UupgbyJiq6Aw9aUI9rm6kXCltcBRlv%2FdDMFxt99agNdK3vi4%2B%2BQZuAKhCALjB2TN2EZA6efhNOAW4y0tORrNVnoWdOmNHXFEh6JlKS6dv1DKPdRav05xoccU4gpLASbjizN2RzhSurAgyato%2FfK2fpMWm1sj1O%2B1rZs8K11fG0Qkr5f7N9emLDENT4xJv%2Fi7m9IKd%2BLf7xHZTfWk8xLaHes4ncg2hO%2FZeMXjQcjum7UQHyWcS2Rfhb7kzA74fxHtjfSW3MtjSU31Hi2NdJVJCNHs8ZX%2FIFZmPJH9gq5hWj3NXAFjpevz5zz0ZWovkGKq7D00RRfp9GqIK%2BYSzYNJ8XqJACflqgQtDzS9ymXvbMK2FPCzNG1667r%2FR6%2BJPSUupGrPwkRw%2B4QmxApO6HX1NzegKE8%2BlUXvyAMGEuyVoVj4uokf5IP1P9%2FgWAlxSAbfRwGhB5VfG9a8rYERYnO8IDU%2B9Tta2MlusSEp5siXPq2PR8QBMqoC%2F8zM7gM7BwND6ZILuXNiqaQbGD3XJ96IYGzAutb1mVofPzyjcBzolV2T88xKwv0FGrE%2FpJPzJEPMl5ouZ6oBxxTDiQO10x53lFQJ0w3pvqwQY6pQGoJi3k6bpfI1as9VEqVsSHHieB00CrBHFBlN%2BX7K7zlpMS12ZGN%2BbSoHkhsk2zLCzRsmc%2BShOln5hQBKyLJfDHn01tz0wDxXnXFyrkzVM4kWjVvA4yMXdN3D%2B%2BuSX%2FcY5bQPijbW78%2FRKP82bnDhqn4yygMiBgsM9j%2BIKWOWaw6mkZpNdTRDzFmRJDEzAj94RO4QEUTNpEd%2B8%2FvOHlpFzkvooAGrg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2031%20May%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=dfd3b7d0da55bbf06e94ce007afe3359917a4520b735281596a1558d24e01cfc"
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991109254/1748621865779/iSPzvhCEiJ-camera-video-segment-1748621868279.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAJAEKMH6%2F20250531%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250531T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCuUK7w6mmYqW2%2FcaY5Zsk%2FUjw3OSnuIDDgNxk%2FJCULYAIgfudyuGCLrAYHoTeAmUhjEBGJ6BPZjUG13vKGt9tCUXQq5gMItv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDF02x%2Bud36UupgbyJiq6Aw9aUI9rm6kXCltcBRlv%2FdDMFxt99agNdK3vi4%2B%2BQZuAKhCALjB2TN2EZA6efhNOAW4y0tORrNVnoWdOmNHXFEh6JlKS6dv1DKPdRav05xoccU4gpLASbjizN2RzhSurAgyato%2FfK2fpMWm1sj1O%2B1rZs8K11fG0Qkr5f7N9emLDENT4xJv%2Fi7m9IKd%2BLf7xHZTfWk8xLaHes4ncg2hO%2FZeMXjQcjum7UQHyWcS2Rfhb7kzA74fxHtjfSW3MtjSU31Hi2NdJVJCNHs8ZX%2FIFZmPJH9gq5hWj3NXAFjpevz5zz0ZWovkGKq7D00RRfp9GqIK%2BYSzYNJ8XqJACflqgQtDzS9ymXvbMK2FPCzNG1667r%2FR6%2BJPSUupGrPwkRw%2B4QmxApO6HX1NzegKE8%2BlUXvyAMGEuyVoVj4uokf5IP1P9%2FgWAlxSAbfRwGhB5VfG9a8rYERYnO8IDU%2B9Tta2MlusSEp5siXPq2PR8QBMqoC%2F8zM7gM7BwND6ZILuXNiqaQbGD3XJ96IYGzAutb1mVofPzyjcBzolV2T88xKwv0FGrE%2FpJPzJEPMl5ouZ6oBxxTDiQO10x53lFQJ0w3pvqwQY6pQGoJi3k6bpfI1as9VEqVsSHHieB00CrBHFBlN%2BX7K7zlpMS12ZGN%2BbSoHkhsk2zLCzRsmc%2BShOln5hQBKyLJfDHn01tz0wDxXnXFyrkzVM4kWjVvA4yMXdN3D%2B%2BuSX%2FcY5bQPijbW78%2FRKP82bnDhqn4yygMiBgsM9j%2BIKWOWaw6mkZpNdTRDzFmRJDEzAj94RO4QEUTNpEd%2B8%2FvOHlpFzkvooAGrg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2031%20May%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=5adf96a41bda7b124603ddfcde2da8358ede071d9e8f4d8083483fdc067e934c"
$ws.Range("L3").Style = "Normal"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "281474989902631-1748616501806"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "Forward Collision Warning"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2025-05-30T08:48:21.806"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "281474989902631"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "151"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "51860288"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "FERNANDO GARCÍA"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 20.67714335
$ws.Range("I4").Value = -103.01918147
$ws.Range("J4").Value = 0
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474989902631/1748616495306/H9b4wJS3sd-camera-video-segment-driver-1748616500306.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAJAEKMH6%2F20250531%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250531T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCuUK7w6mmYqW2%2FcaY5Zsk%2FUjw3OSnuIDDgNxk%2FJCULYAIgfudyuGCLrAYHoTeAmUhjEBGJ6BPZjUG13vKGt9tCUXQq5gMItv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDF02x%2Bud36UupgbyJiq6Aw9aUI9rm6kXCltcBRlv%2FdDMFxt99agNdK3vi4%2B%2BQZuAKhCALjB2TN2EZA6efhNOAW4y0tORrNVnoWdOmNHXFEh6JlKS6dv1DKPdRav05xoccU4gpLASbjizN2RzhSurAgyato%2FfK2fpMWm1sj1O%2B1rZs8K11fG0Qkr5f7N9emLDENT4xJv%2Fi7m9IKd%2BLf7xHZTfWk8xLaHes4ncg2hO%2FZeMXjQcjum7UQHyWcS2Rfhb7kzA74fxHtjfSW3MtjSU31Hi2NdJVJCNHs8ZX%2FIFZmPJH9gq5hWj3NXAFjpevz5zz0ZWovkGKq7D00RRfp9GqIK%2BYSzYNJ8XqJACflqgQtDzS9ymXvbMK2FPCzNG1667r%2FR6%2BJPSUupGrPwkRw%2B4QmxApO6HX1NzegKE8%2BlUXvyAMGEuyVoVj4uokf5IP1P9%2FgWAlxSAbfRwGhB5VfG9a8rYERYnO8IDU%2B9Tta2MlusSEp5siXPq2PR8QBMqoC%2F8zM7gM7BwND6ZILuXNiqaQbGD3XJ96IYGzAutb1mVofPzyjcBzolV2T88xKwv0FGrE%2FpJPzJEPMl5ouZ6oBxxTDiQO10x53lFQJ0w3pvqwQY6pQGoJi3k6bpfI1as9VEqVsSHHieB00CrBHFBlN%2BX7K7zlpMS12ZGN%2BbSoHkhsk2zLCzRsmc%2BShOln5hQBKyLJfDHn01tz0wDxXnXFyrkzVM4kWjVvA4yMXdN3D%2B%2BuSX%2FcY5bQPijbW78%2FRKP82bnDhqn4yygMiBgsM9j%2BIKWOWaw6mkZpNdTRDzFmRJDEzAj94RO4QEUTNpEd%2B8%2FvOHlpFzkvooAGrg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2031%20May%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=49d3de5cd54c387f1d8617a821b40026e0630fbe9836209f7307bff202e25a9c"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474989902631/1748616495306/tLMzaoUUhQ-camera-video-segment-1748616500306.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAJAEKMH6%2F20250531%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250531T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCuUK7w6mmYqW2%2FcaY5Zsk%2FUjw3OSnuIDDgNxk%2FJCULYAIgfudyuGCLrAYHoTeAmUhjEBGJ6BPZjUG13vKGt9tCUXQq5gMItv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDF02x%2Bud36UupgbyJiq6Aw9aUI9rm6kXCltcBRlv%2FdDMFxt99agNdK3vi4%2B%2BQZuAKhCALjB2TN2EZA6efhNOAW4y0tORrNVnoWdOmNHXFEh6JlKS6dv1DKPdRav05xoccU4gpLASbjizN2RzhSurAgyato%2FfK2fpMWm1sj1O%2B1rZs8K11fG0Qkr5f7N9emLDENT4xJv%2Fi7m9IKd%2BLf7xHZTfWk8xLaHes4ncg2hO%2FZeMXjQcjum7UQHyWcS2Rfhb7kzA74fxHtjfSW3MtjSU31Hi2NdJVJCNHs8ZX%2FIFZmPJH9gq5hWj3NXAFjpevz5zz0ZWovkGKq7D00RRfp9GqIK%2BYSzYNJ8XqJACflqgQtDzS9ymXvbMK2FPCzNG1667r%2FR6%2BJPSUupGrPwkRw%2B4QmxApO6HX1NzegKE8%2BlUXvyAMGEuyVoVj4uokf5IP1P9%2FgWAlxSAbfRwGhB5VfG9a8rYERYnO8IDU%2B9Tta2MlusSEp5siXPq2PR8QBMqoC%2F8zM7gM7BwND6ZILuXNiqaQbGD3XJ96IYGzAutb1mVofPzyjcBzolV2T88xKwv0FGrE%2FpJPzJEPMl5ouZ6oBxxTDiQO10x53lFQJ0w3pvqwQY6pQGoJi3k6bpfI1as9VEqVsSHHieB00CrBHFBlN%2BX7K7zlpMS12ZGN%2BbSoHkhsk2zLCzRsmc%2BShOln5hQBKyLJfDHn01tz0wDxXnXFyrkzVM4kWjVvA4yMXdN3D%2B%2BuSX%2FcY5bQPijbW78%2FRKP82bnDhqn4yygMiBgsM9j%2BIKWOWaw6mkZpNdTRDzFmRJDEzAj94RO4QEUTNpEd%2B8%2FvOHlpFzkvooAGrg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2031%20May%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=00813bf30d371346184fd5b6c0d565e747863fdd1f0f68a336e47a1902830ee6"
$ws.Range("L4").Style = "Normal"
